$d = $word.ActiveDocument

# 1. Append " for enterprise pro" after "Meeting Minutes" in the title paragraph,
#    matching the existing run formatting (bold, underline, size 40/20pt).
$titlePara = $d.Paragraphs(1)
$endRange = $titlePara.Range
$endRange.Collapse(0)
$endRange.Font.Bold = $true
$endRange.Font.Underline = 1
$endRange.Font.Size = 20
$endRange.InsertAfter(" for enterprise pro")

# 2. Merge the "Plan for the next week" + ":" runs into a single run (no visible text change).
$d.Content.Find.Execute(" Plan for the next week:", $false, $false, $false, $false, $false, $true, 1, $false, " Plan for the next week:", 2)
